$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra empty row (row 2) that sits between the header and the
# data, shifting everything below it up by one row.
$ws.Rows("2").Delete()

# Update the selection to match the post-edit state captured in the diff.
$ws.Range("F8").Select()
